{"js": "// Replace each three-digit-divided-by-one-digit expression in the table\n// with its updated value, per the commit's regenerated problem set.\nconst replacements = [\n  [\"733\u00f78=\", \"591\u00f76=\"],\n  [\"911\u00f76=\", \"174\u00f72=\"],\n  [\"523\u00f78=\", \"312\u00f75=\"],\n  [\"228\u00f77=\", \"341\u00f76=\"],\n  [\"957\u00f78=\", \"578\u00f74=\"],\n  [\"352\u00f78=\", \"557\u00f77=\"],\n  [\"536\u00f72=\", \"997\u00f72=\"],\n  [\"229\u00f78=\", \"666\u00f72=\"],\n  [\"913\u00f72=\", \"876\u00f78=\"],\n  [\"424\u00f76=\", \"715\u00f76=\"],\n  [\"489\u00f79=\", \"558\u00f77=\"],\n  [\"287\u00f75=\", \"255\u00f78=\"],\n  [\"694\u00f78=\", \"709\u00f74=\"],\n  [\"182\u00f77=\", \"343\u00f76=\"],\n  [\"533\u00f75=\", \"134\u00f78=\"],\n  [\"534\u00f72=\", \"277\u00f79=\"],\n  [\"391\u00f75=\", \"683\u00f79=\"],\n  [\"909\u00f75=\", \"129\u00f74=\"],\n  [\"178\u00f77=\", \"141\u00f73=\"],\n  [\"699\u00f74=\", \"307\u00f76=\"],\n  [\"669\u00f74=\", \"159\u00f75=\"],\n  [\"925\u00f77=\", \"876\u00f72=\"],\n  [\"978\u00f74=\", \"129\u00f72=\"],\n  [\"817\u00f79=\", \"796\u00f73=\"],\n  [\"192\u00f78=\", \"661\u00f74=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update master to output generated at 4250d90\n# Replace each three-digit / one-digit division expression in the practice\n# table with its regenerated value (old => new, one-to-one, order-independent).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"733\u00f78=\", \"591\u00f76=\"),\n    @(\"911\u00f76=\", \"174\u00f72=\"),\n    @(\"523\u00f78=\", \"312\u00f75=\"),\n    @(\"228\u00f77=\", \"341\u00f76=\"),\n    @(\"957\u00f78=\", \"578\u00f74=\"),\n    @(\"352\u00f78=\", \"557\u00f77=\"),\n    @(\"536\u00f72=\", \"997\u00f72=\"),\n    @(\"229\u00f78=\", \"666\u00f72=\"),\n    @(\"913\u00f72=\", \"876\u00f78=\"),\n    @(\"424\u00f76=\", \"715\u00f76=\"),\n    @(\"489\u00f79=\", \"558\u00f77=\"),\n    @(\"287\u00f75=\", \"255\u00f78=\"),\n    @(\"694\u00f78=\", \"709\u00f74=\"),\n    @(\"182\u00f77=\", \"343\u00f76=\"),\n    @(\"533\u00f75=\", \"134\u00f78=\"),\n    @(\"534\u00f72=\", \"277\u00f79=\"),\n    @(\"391\u00f75=\", \"683\u00f79=\"),\n    @(\"909\u00f75=\", \"129\u00f74=\"),\n    @(\"178\u00f77=\", \"141\u00f73=\"),\n    @(\"699\u00f74=\", \"307\u00f76=\"),\n    @(\"669\u00f74=\", \"159\u00f75=\"),\n    @(\"925\u00f77=\", \"876\u00f72=\"),\n    @(\"978\u00f74=\", \"129\u00f72=\"),\n    @(\"817\u00f79=\", \"796\u00f73=\"),\n    @(\"192\u00f78=\", \"661\u00f74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
